# Apply trade-close update + new open trade to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet - headline stats refreshed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.95   # Current Capital
$summary.Range("B4").Value = 0.75      # Total P&L $
$summary.Range("B5").Value = 0.18      # Total P&L %
$summary.Range("B6").Value = 85        # Total Trades
$summary.Range("B7").Value = 40        # Winning Trades
$summary.Range("B9").Value = 47.06     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5) refreshed
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.95     # Capital
$status.Range("D5").Value = 52         # Trades
$status.Range("E5").Value = 0.64       # P&L $
$status.Range("F5").Value = 0.95       # P&L %
$status.Range("G5").Value = 50         # Win Rate %

# ---------------------------------------------------------------------------
# 3) All Trades sheet - close out trade #85 (row 86) and append trade #118
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G86").Value = 0.263158
$allTrades.Range("H86").Value = "CLOSED"
$allTrades.Range("I86").Value = 19.6172
$allTrades.Range("J86").Value = 0.04
$allTrades.Range("K86").Value = 100.95
$allTrades.Range("L86").Value = "early_exit"
$allTrades.Range("M86").Value = 0.13

$allTrades.Range("A119").Value = 118
$allTrades.Range("B119").NumberFormat = "@"
$allTrades.Range("B119").Value = "2026-02-17"
$allTrades.Range("B119").Style = "Normal"
$allTrades.Range("C119").Value = "21:10:24"
$allTrades.Range("D119").Value = "MarketMaking"
$allTrades.Range("E119").Value = "DOWN"
$allTrades.Range("F119").Value = 0.22
$allTrades.Range("H119").Value = "OPEN"
$allTrades.Range("I119").Value = 0
$allTrades.Range("J119").Value = 0
$allTrades.Range("K119").Value = 100.9114872031006
$allTrades.Range("M119").Value = 0
$allTrades.Range("N119").Value = 0
$allTrades.Range("O119").Value = 0
$allTrades.Range("P119").Value = 0.6
$allTrades.Range("Q119").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# 4) MarketMaking sheet - close out trade #85 (row 53) and append trade #118
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G53").Value = 0.263158
$mm.Range("H53").Value = "CLOSED"
$mm.Range("I53").Value = 19.6172
$mm.Range("J53").Value = 0.04
$mm.Range("K53").Value = 100.95
$mm.Range("P53").Value = "early_exit"
$mm.Range("Q53").Value = 0.13

$mm.Range("A86").Value = 118
$mm.Range("B86").NumberFormat = "@"
$mm.Range("B86").Value = "2026-02-17"
$mm.Range("B86").Style = "Normal"
$mm.Range("C86").Value = "21:10:24"
$mm.Range("D86").Value = "MarketMaking"
$mm.Range("E86").Value = "DOWN"
$mm.Range("F86").Value = 0.22
$mm.Range("H86").Value = "OPEN"
$mm.Range("I86").Value = 0
$mm.Range("J86").Value = 0
$mm.Range("K86").Value = 100.9114872031006
$mm.Range("L86").Value = 0
$mm.Range("M86").Value = 0
$mm.Range("N86").Value = 0.6
$mm.Range("O86").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q86").Value = 0
